$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "AAPL"
$ws.Range("B3").Value = 231.5899963378906

$ws.Range("A4").Value = "MSFT"
$ws.Range("B4").Value = 520.1699829101562

$ws.Range("A5").Value = "GOOGL"
$ws.Range("B5").Value = 203.8999938964844
